{"js": "const doc = context.document;\nconst body = doc.body;\n\n// 1. Remove the \"_GoBack\" bookmark from its current position (right before\n//    the run \"this information and keep the screen blank.\").\ndoc.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2. Locate the paragraph that holds the embedded ActiveX/OLE control. It\n//    has no visible text (the control isn't represented as plain text) and\n//    is immediately followed by the document's final, empty paragraph.\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nparas.items.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\nconst count = paras.items.length;\nlet targetIndex = -1;\nif (count >= 2 &&\n    paras.items[count - 1].text === \"\" &&\n    paras.items[count - 2].text === \"\") {\n  targetIndex = count - 2;\n}\n\nif (targetIndex !== -1) {\n  paras.items[targetIndex].delete();\n  await context.sync();\n}\n\n// 3. Re-add the \"_GoBack\" bookmark to the (now last) trailing empty\n//    paragraph.\nconst paras2 = body.paragraphs;\nparas2.load(\"items\");\nawait context.sync();\n\nconst last = paras2.items[paras2.items.length - 1];\nconst lastRange = last.getRange(\"Whole\");\nlastRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Remove the \"_GoBack\" bookmark from its current position (right before the\n#    run \"this information and keep the screen blank.\").\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2. Locate and delete the paragraph that holds the embedded ActiveX/OLE\n#    control (it is shaded F2F3F5, a color unique to that paragraph) \u2014 it\n#    sits right before the final, empty paragraph.\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Shading.BackgroundPatternColor -eq 16118770) {\n        $target = $p\n    }\n}\nif ($target -ne $null) {\n    $target.Range.Delete()\n}\n\n# 3. Re-add the \"_GoBack\" bookmark to the (now last) trailing empty paragraph.\n$last = $d.Paragraphs.Item($d.Paragraphs.Count)\n$d.Bookmarks.Add(\"_GoBack\", $last.Range)\n"}
